# Apply updated market/profit figures to the Leve profit sheets.
# Values correspond to refreshed currentAveragePrice / LevePrice / LeveProfit
# figures pulled by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 64999
$ws.Range("J108").Value = 64999
$ws.Range("L108").Value = 64999
$ws.Range("N108").Value = -72679
$ws.Range("H110").Value = 38500
$ws.Range("J110").Value = 38500
$ws.Range("L110").Value = 38500
$ws.Range("N110").Value = -46680
$ws.Range("H112").Value = 2189.6792
$ws.Range("I112").Value = 949.6
$ws.Range("J112").Value = 2318.8542
$ws.Range("K112").Value = 2848.8
$ws.Range("L112").Value = 6956.562600000001
$ws.Range("M112").Value = -1740.8
$ws.Range("N112").Value = -9172.562600000001
$ws.Range("H116").Value = 5294862
$ws.Range("I116").Value = 9262834
$ws.Range("J116").Value = 4233.6665
$ws.Range("K116").Value = 9262834
$ws.Range("L116").Value = 4233.6665
$ws.Range("M116").Value = -9259392
$ws.Range("N116").Value = -11117.6665
$ws.Range("H124").Value = 200000
$ws.Range("I124").Value = 200000
$ws.Range("K124").Value = 200000
$ws.Range("M124").Value = -195090
$ws.Range("H132").Value = 260346.2
$ws.Range("I132").Value = 281380.62
$ws.Range("K132").Value = 844141.86
$ws.Range("M132").Value = -841611.86

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 715938.5600000001
$ws.Range("I31").Value = 715938.5600000001
$ws.Range("K31").Value = 715938.5600000001
$ws.Range("M31").Value = -715644.5600000001
$ws.Range("H32").Value = 1955663.4
$ws.Range("I32").Value = 4760.109
$ws.Range("J32").Value = 10113986
$ws.Range("K32").Value = 4760.109
$ws.Range("L32").Value = 10113986
$ws.Range("M32").Value = -4473.109
$ws.Range("N32").Value = -10114560
$ws.Range("H61").Value = 4679.2334
$ws.Range("I61").Value = 4589.55
$ws.Range("J61").Value = 4724.075
$ws.Range("K61").Value = 4589.55
$ws.Range("L61").Value = 4724.075
$ws.Range("M61").Value = -4377.55
$ws.Range("N61").Value = -5148.075
$ws.Range("H136").Value = 4679.2334
$ws.Range("I136").Value = 4589.55
$ws.Range("J136").Value = 4724.075
$ws.Range("K136").Value = 13768.65
$ws.Range("L136").Value = 14172.225
$ws.Range("M136").Value = -11218.65
$ws.Range("N136").Value = -19272.225

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6494537.5
$ws.Range("I20").Value = 9524795
$ws.Range("J20").Value = 1128.5714
$ws.Range("K20").Value = 9524795
$ws.Range("L20").Value = 1128.5714
$ws.Range("M20").Value = -9524548
$ws.Range("N20").Value = -1622.5714
$ws.Range("H22").Value = 1958.0834
$ws.Range("I22").Value = 336.5
$ws.Range("K22").Value = 336.5
$ws.Range("M22").Value = -163.5
$ws.Range("H128").Value = 4626
$ws.Range("I128").Value = 4626
$ws.Range("K128").Value = 13878
$ws.Range("M128").Value = -11388
$ws.Range("H134").Value = 1116609
$ws.Range("I134").Value = 1320620.6
$ws.Range("J134").Value = 9117.429
$ws.Range("K134").Value = 3961861.8
$ws.Range("L134").Value = 27352.287
$ws.Range("M134").Value = -3959326.8
$ws.Range("N134").Value = -32422.287

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 659.3077
$ws.Range("I7").Value = 94.111115
$ws.Range("K7").Value = 94.111115
$ws.Range("M7").Value = 18.888885
$ws.Range("H22").Value = 812394
$ws.Range("J22").Value = 1306
$ws.Range("L22").Value = 1306
$ws.Range("N22").Value = -2006
$ws.Range("H58").Value = 27787564
$ws.Range("I58").Value = 31257338
$ws.Range("J58").Value = 29374.75
$ws.Range("K58").Value = 31257338
$ws.Range("L58").Value = 29374.75
$ws.Range("M58").Value = -31257135
$ws.Range("N58").Value = -29780.75
$ws.Range("H99").Value = 5559027
$ws.Range("I99").Value = 10104053
$ws.Range("K99").Value = 10104053
$ws.Range("M99").Value = -10102555
$ws.Range("H122").Value = 10374.777
$ws.Range("I122").Value = 2953.3333
$ws.Range("J122").Value = 14085.5
$ws.Range("K122").Value = 8859.999899999999
$ws.Range("L122").Value = 42256.5
$ws.Range("M122").Value = -6409.999899999999
$ws.Range("N122").Value = -47156.5
$ws.Range("H126").Value = 5559027
$ws.Range("I126").Value = 10104053
$ws.Range("K126").Value = 30312159
$ws.Range("M126").Value = -30309689
$ws.Range("H136").Value = 27787564
$ws.Range("I136").Value = 31257338
$ws.Range("J136").Value = 29374.75
$ws.Range("K136").Value = 93772014
$ws.Range("L136").Value = 88124.25
$ws.Range("M136").Value = -93769464
$ws.Range("N136").Value = -93224.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2737.4211
$ws.Range("I117").Value = 778.75
$ws.Range("K117").Value = 2336.25
$ws.Range("M117").Value = 1105.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7200.5
$ws.Range("I70").Value = 6669.5
$ws.Range("J70").Value = 7997
$ws.Range("K70").Value = 6669.5
$ws.Range("L70").Value = 7997
$ws.Range("M70").Value = -6399.5
$ws.Range("N70").Value = -8537
$ws.Range("H73").Value = 7200.5
$ws.Range("I73").Value = 6669.5
$ws.Range("J73").Value = 7997
$ws.Range("K73").Value = 6669.5
$ws.Range("L73").Value = 7997
$ws.Range("M73").Value = -5733.5
$ws.Range("N73").Value = -9869
$ws.Range("H80").Value = 8474.6
$ws.Range("J80").Value = 9231
$ws.Range("L80").Value = 9231
$ws.Range("N80").Value = -11227
$ws.Range("H83").Value = 8474.6
$ws.Range("J83").Value = 9231
$ws.Range("L83").Value = 46155
$ws.Range("N83").Value = -56139
$ws.Range("H97").Value = 2334.3333
$ws.Range("I97").Value = 2014.1765
$ws.Range("K97").Value = 2014.1765
$ws.Range("M97").Value = -1518.1765
$ws.Range("H102").Value = 942563.2
$ws.Range("I102").Value = 1555517.5
$ws.Range("K102").Value = 1555517.5
$ws.Range("M102").Value = -1553895.5
$ws.Range("H132").Value = 5259.0527
$ws.Range("I132").Value = 5363.3784
$ws.Range("J132").Value = 1399
$ws.Range("K132").Value = 16090.1352
$ws.Range("L132").Value = 4197
$ws.Range("M132").Value = -13560.1352
$ws.Range("N132").Value = -9257
$ws.Range("H134").Value = 71150
$ws.Range("J134").Value = 71150
$ws.Range("L134").Value = 213450
$ws.Range("N134").Value = -218520

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 978.1667
$ws.Range("I22").Value = 969.75
$ws.Range("K22").Value = 969.75
$ws.Range("M22").Value = -674.75
$ws.Range("H27").Value = 978.1667
$ws.Range("I27").Value = 969.75
$ws.Range("K27").Value = 969.75
$ws.Range("M27").Value = -862.75
$ws.Range("H93").Value = 3632.7778
$ws.Range("I93").Value = 3724.375
$ws.Range("J93").Value = 2900
$ws.Range("K93").Value = 3724.375
$ws.Range("L93").Value = 2900
$ws.Range("M93").Value = -2476.375
$ws.Range("N93").Value = -5396
$ws.Range("H132").Value = 2464.5
$ws.Range("I132").Value = 2464.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7393.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4863.5
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 35722708
$ws.Range("I136").Value = 55562664
$ws.Range("K136").Value = 166687992
$ws.Range("M136").Value = -166685442

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3971.5715
$ws.Range("I96").Value = 4050.25
$ws.Range("J96").Value = 3866.6667
$ws.Range("K96").Value = 4050.25
$ws.Range("L96").Value = 3866.6667
$ws.Range("M96").Value = -2677.25
$ws.Range("N96").Value = -6612.6667
$ws.Range("H132").Value = 9558.243
$ws.Range("I132").Value = 5670.433
$ws.Range("K132").Value = 17011.299
$ws.Range("M132").Value = -14481.299
$ws.Range("H136").Value = 21753962
$ws.Range("I136").Value = 23820780
$ws.Range("K136").Value = 71462340
$ws.Range("M136").Value = -71459790
